# Update cryptocurrency price/volume data in the worksheet
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text looks like a plain number need an explicit
# text number-format first, otherwise Excel would coerce the literal
# "27.282.86"-style price strings into numeric values.
$numericLookingCells = @('D4','D5','D7','D8','D9','D10','D11','D12','D13','D14','D15','D17','D18','D19','D20','D23','D26','D27','D28','D29','D31','D32','D33','D34','D35','D36','D37','D38','D39','D40','D41','D43','D44','D45','D46','D47','D48','D49','D51')
foreach ($addr in $numericLookingCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range('D2').Value2 = '27.282.86'
$ws.Range('E2').Value2 = '  -3.12%  '
$ws.Range('D3').Value2 = '1.731.85'
$ws.Range('E3').Value2 = '  -3.88%  '
$ws.Range('D4').Value2 = '1.005'
$ws.Range('E4').Value2 = '  +0.11%  '
$ws.Range('D5').Value2 = '322.57'
$ws.Range('E5').Value2 = '  -4.88%  '
$ws.Range('E6').Value2 = '  +0.07%  '
$ws.Range('D7').Value2 = '0.4243'
$ws.Range('E7').Value2 = '  -10.92%  '
$ws.Range('D8').Value2 = '0.3579'
$ws.Range('E8').Value2 = '  -4.50%  '
$ws.Range('D9').Value2 = '44.80'
$ws.Range('E9').Value2 = '  -1.36%  '
$ws.Range('D10').Value2 = '0.07456'
$ws.Range('D11').Value2 = '1.110'
$ws.Range('E11').Value2 = '  -3.50%  '
$ws.Range('D12').Value2 = '1.003'
$ws.Range('E12').Value2 = '  -0.03%  '
$ws.Range('D13').Value2 = '21.38'
$ws.Range('E13').Value2 = '  -5.74%  '
$ws.Range('D14').Value2 = '6.051'
$ws.Range('E14').Value2 = '  -5.46%  '
$ws.Range('D15').Value2 = '7.139'
$ws.Range('E15').Value2 = '  -3.73%  '
$ws.Range('D16').Value2 = '1.731.77'
$ws.Range('E16').Value2 = '  -3.74%  '
$ws.Range('D17').Value2 = '0.00001062'
$ws.Range('E17').Value2 = '  -3.29%  '
$ws.Range('D18').Value2 = '86.76'
$ws.Range('E18').Value2 = '  +5.00%  '
$ws.Range('D19').Value2 = '0.05991'
$ws.Range('E19').Value2 = '  -11.31%  '
$ws.Range('D20').Value2 = '1.002'
$ws.Range('E20').Value2 = '  +0.11%  '
$ws.Range('E21').Value2 = '  -4.60%  '
$ws.Range('E22').Value2 = '  -6.21%  '
$ws.Range('D23').Value2 = '0.5234'
$ws.Range('E23').Value2 = '  -5.44%  '
$ws.Range('D24').Value2 = '27.305.94'
$ws.Range('E24').Value2 = '  -3.00%  '
$ws.Range('E25').Value2 = '  -5.35%  '
$ws.Range('D26').Value2 = '2.389'
$ws.Range('E26').Value2 = '  -0.86%  '
$ws.Range('D27').Value2 = '20.08'
$ws.Range('E27').Value2 = '  -4.08%  '
$ws.Range('D28').Value2 = '2.346'
$ws.Range('E28').Value2 = '  -2.34%  '
$ws.Range('D29').Value2 = '149.73'
$ws.Range('E29').Value2 = '  -1.00%  '
$ws.Range('D30').Value2 = '1.928.85'
$ws.Range('E30').Value2 = '  -3.77%  '
$ws.Range('D31').Value2 = '1.249'
$ws.Range('E31').Value2 = '  -1.30%  '
$ws.Range('D32').Value2 = '126.02'
$ws.Range('E32').Value2 = '  -6.39%  '
$ws.Range('D33').Value2 = '3.695'
$ws.Range('E33').Value2 = '  -8.69%  '
$ws.Range('D34').Value2 = '0.09060'
$ws.Range('E34').Value2 = '  -6.70%  '
$ws.Range('D35').Value2 = '5.576'
$ws.Range('E35').Value2 = '  -6.41%  '
$ws.Range('D36').Value2 = '12.53'
$ws.Range('E36').Value2 = '  +2.32%  '
$ws.Range('D37').Value2 = '0.2150'
$ws.Range('E37').Value2 = '  -3.56%  '
$ws.Range('D38').Value2 = '0.02276'
$ws.Range('E38').Value2 = '  -4.86%  '
$ws.Range('D39').Value2 = '0.06122'
$ws.Range('E39').Value2 = '  -4.03%  '
$ws.Range('D40').Value2 = '5.023'
$ws.Range('E40').Value2 = '  -5.06%  '
$ws.Range('D41').Value2 = '0.6367'
$ws.Range('E41').Value2 = '  -5.62%  '
$ws.Range('E42').Value2 = '  -4.68%  '
$ws.Range('D43').Value2 = '1.002'
$ws.Range('E43').Value2 = '  +0.10%  '
$ws.Range('D44').Value2 = '1.405'
$ws.Range('E44').Value2 = '  -5.28%  '
$ws.Range('D45').Value2 = '7.855'
$ws.Range('E45').Value2 = '  -3.09%  '
$ws.Range('D46').Value2 = '13.48'
$ws.Range('E46').Value2 = '  -5.31%  '
$ws.Range('D47').Value2 = '3.732'
$ws.Range('E47').Value2 = '  -3.41%  '
$ws.Range('D48').Value2 = '0.5815'
$ws.Range('E48').Value2 = '  -6.02%  '
$ws.Range('D49').Value2 = '124.19'
$ws.Range('E49').Value2 = '  -5.10%  '
$ws.Range('E50').Value2 = '  -6.47%  '
$ws.Range('D51').Value2 = '0.06817'
$ws.Range('E51').Value2 = '  -4.37%  '

# Restore the default cell style so only the values changed.
foreach ($addr in $numericLookingCells) {
    $ws.Range($addr).Style = "Normal"
}
